{"js": "// Replace the 31 numbered paragraphs (\"1\"..\"31\") with a single paragraph\n// containing two lines of Lorem-ipsum text joined by a manual line break\n// (w:br), keeping the paragraph's line spacing (360/auto) but dropping the\n// w:after=\"0\" and w:jc=\"left\" paragraph-formatting overrides.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Drop every paragraph except the first one (its pPr/run will be\n// overwritten below via insertOoxml, so which one survives doesn't matter\n// as long as exactly one remains for the sectPr to stay attached to the\n// body correctly).\nfor (let i = paragraphs.items.length - 1; i >= 1; i--) {\n  paragraphs.items[i].delete();\n}\nawait context.sync();\n\nconst line1 =\n  \"Lorem ipsum odor amet, consectetuer adipiscing elit. Congue vel parturient sapien volutpat porttitor malesuada mus. Volutpat sociosqu nisi cubilia himenaeos sed in nisl leo. Dis venenatis ullamcorper pharetra; penatibus blandit arcu justo dignissim nullam. Dolor a sodales, nostra lacinia nascetur faucibus. Sodales volutpat mattis suscipit morbi\";\nconst line2 =\n  \"elementum sapien convallis nec egestas. Dignissim lacinia dolor placerat nulla porta natoque fames, sem non. Venenatis facilisi dapibus pellentesque netus etiam id blandit. Porttitor integer nec urna posuere rhoncus faucibus.\";\n\nfunction escapeXml(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\nconst target = remaining.items[0];\nconst targetRange = target.getRange();\n\n// Use insertOoxml (replace) so we can set the paragraph mark's pPr and the\n// run's two <w:t> runs split by <w:br/> exactly, instead of accumulating\n// separate runs/paragraphs the way insertText/insertBreak would.\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>${escapeXml(line1)}</w:t>\n              <w:br/>\n              <w:t>${escapeXml(line2)}</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntargetRange.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the 31 numbered paragraphs (\"1\"..\"31\") with a single paragraph\n# containing two lines of Lorem-ipsum text joined by a manual line break\n# (w:br), keeping the paragraph's line spacing (360/auto) but dropping the\n# w:after=\"0\" and w:jc=\"left\" paragraph-formatting overrides.\n\n$d = $word.ActiveDocument\n\n# Remove every paragraph after the first one.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 2; $i--) {\n    $d.Paragraphs($i).Range.Delete()\n}\n\n$line1 = \"Lorem ipsum odor amet, consectetuer adipiscing elit. Congue vel parturient sapien volutpat porttitor malesuada mus. Volutpat sociosqu nisi cubilia himenaeos sed in nisl leo. Dis venenatis ullamcorper pharetra; penatibus blandit arcu justo dignissim nullam. Dolor a sodales, nostra lacinia nascetur faucibus. Sodales volutpat mattis suscipit morbi\"\n$line2 = \"elementum sapien convallis nec egestas. Dignissim lacinia dolor placerat nulla porta natoque fames, sem non. Venenatis facilisi dapibus pellentesque netus etiam id blandit. Porttitor integer nec urna posuere rhoncus faucibus.\"\n\n$line1x = $line1 -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n$line2x = $line2 -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n\n# Build a minimal single-part OOXML package describing the desired\n# paragraph (pPr with only line spacing + one run holding two <w:t> split\n# by <w:br/>), then use InsertXML to overwrite the whole document body\n# content with it (InsertXML replaces the contents of the range it is\n# called on -- using $d.Content, the whole-body range, avoids leaving a\n# stray extra paragraph behind).\n$ooxml = @\"\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:spacing w:line=\"360\" w:lineRule=\"auto\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>$line1x</w:t>\n              <w:br/>\n              <w:t>$line2x</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$target = $d.Content\n$null = $target.InsertXML($ooxml)\n"}
